$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything down by two rows: old row 1 (headers) -> row 3,
# old rows 3..14 -> rows 5..16, old row 16 -> row 18. Keeps values/styles intact.
$ws.Rows("1:2").Insert()

# New leading label cell for the (now) sub-header row, and rename the
# "base_clust_N" column headers to "base_subtype_N".
$ws.Range("A3").Value2 = "Comorbidities"
$ws.Range("C3").Value2 = "base_subtype_1"
$ws.Range("E3").Value2 = "base_subtype_2"
$ws.Range("G3").Value2 = "base_subtype_3"
$ws.Range("I3").Value2 = "base_subtype_4"
$ws.Range("K3").Value2 = "base_subtype_5"

# New title row, bold + centered, merged across B1:K1.
$ws.Range("B1:K1").Font.Bold = $true
$ws.Range("B1:K1").HorizontalAlignment = -4108
$ws.Range("B1").Value2 = "Supplementary table 1 : comparison of outcomes and main clinical features in the original data ('base') and the cross-validation ('partition') over 10 iterations"
$ws.Range("B1:K1").Merge()

# Columns got a bit wider to fit the new "base_subtype_N" labels.
$ws.Columns("A").ColumnWidth = 11.6
$ws.Columns("C").ColumnWidth = 11.6
$ws.Columns("E").ColumnWidth = 11.6
$ws.Columns("G").ColumnWidth = 11.6
$ws.Columns("I").ColumnWidth = 11.6
$ws.Columns("K").ColumnWidth = 11.6

$ws.Range("E8").Select()
